{"js": "const table = context.document.body.tables.getFirst();\nconst updates = [\n  [0, 0, \"30+28=\"],\n  [0, 1, \"80-29=\"],\n  [0, 2, \"55-33=\"],\n  [0, 3, \"14+73=\"],\n  [0, 4, \"29-13=\"],\n  [1, 0, \"30+54=\"],\n  [1, 1, \"91-43=\"],\n  [1, 2, \"45-3=\"],\n  [1, 3, \"76-24=\"],\n  [1, 4, \"98-39=\"],\n  [2, 0, \"46+27=\"],\n  [2, 1, \"59+36=\"],\n  [2, 2, \"67-53=\"],\n  [2, 3, \"51-27=\"],\n  [2, 4, \"5+53=\"],\n  [3, 0, \"38+46=\"],\n  [3, 1, \"56-15=\"],\n  [3, 2, \"24+62=\"],\n  [3, 3, \"77-32=\"],\n  [3, 4, \"93-14=\"],\n  [4, 0, \"50-25=\"],\n  [4, 1, \"96-57=\"],\n  [4, 2, \"82-61=\"],\n  [4, 3, \"84+1=\"],\n  [4, 4, \"89-8=\"],\n  [5, 0, \"64-5=\"],\n  [5, 1, \"98-26=\"],\n  [5, 2, \"44+41=\"],\n  [5, 3, \"89-77=\"],\n  [5, 4, \"53+36=\"],\n  [6, 0, \"73-26=\"],\n  [6, 1, \"20+70=\"],\n  [6, 2, \"45+29=\"],\n  [6, 3, \"68-38=\"],\n  [6, 4, \"77-19=\"],\n  [7, 0, \"35+49=\"],\n  [7, 1, \"11+33=\"],\n  [7, 2, \"2+9=\"],\n  [7, 3, \"48-47=\"],\n  [7, 4, \"12+20=\"],\n  [8, 0, \"90-16=\"],\n  [8, 1, \"61-34=\"],\n  [8, 2, \"84+6=\"],\n  [8, 3, \"23-10=\"],\n  [8, 4, \"11+78=\"],\n  [9, 0, \"49-4=\"],\n  [9, 1, \"42+40=\"],\n  [9, 2, \"84-30=\"],\n  [9, 3, \"96-9=\"],\n  [9, 4, \"29+31=\"],\n  [10, 0, \"49-13=\"],\n  [10, 1, \"52+15=\"],\n  [10, 2, \"67+21=\"],\n  [10, 3, \"25+70=\"],\n  [10, 4, \"43-13=\"],\n  [11, 0, \"25+62=\"],\n  [11, 1, \"86-13=\"],\n  [11, 2, \"52+28=\"],\n  [11, 3, \"77-26=\"],\n  [11, 4, \"77-32=\"],\n  [12, 0, \"14+48=\"],\n  [12, 1, \"21-15=\"],\n  [12, 2, \"6+1=\"],\n  [12, 3, \"10+26=\"],\n  [12, 4, \"45-14=\"],\n  [13, 0, \"70-25=\"],\n  [13, 1, \"2+8=\"],\n  [13, 2, \"44-42=\"],\n  [13, 3, \"10+12=\"],\n  [13, 4, \"91-16=\"],\n  [14, 0, \"48-9=\"],\n  [14, 1, \"14+1=\"],\n  [14, 2, \"87-36=\"],\n  [14, 3, \"85-79=\"],\n  [14, 4, \"34-25=\"],\n  [15, 0, \"99-5=\"],\n  [15, 1, \"19+14=\"],\n  [15, 2, \"47-42=\"],\n  [15, 3, \"23+15=\"],\n  [15, 4, \"13+38=\"],\n  [16, 0, \"66-2=\"],\n  [16, 1, \"12-4=\"],\n  [16, 2, \"91-40=\"],\n  [16, 3, \"17+75=\"],\n  [16, 4, \"87-87=\"],\n  [17, 0, \"10+45=\"],\n  [17, 1, \"90-60=\"],\n  [17, 2, \"1+68=\"],\n  [17, 3, \"58-31=\"],\n  [17, 4, \"99-22=\"],\n  [18, 0, \"49-12=\"],\n  [18, 1, \"64-61=\"],\n  [18, 2, \"61-50=\"],\n  [18, 3, \"44-35=\"],\n  [18, 4, \"35-13=\"],\n  [19, 0, \"82+15=\"],\n  [19, 1, \"79-1=\"],\n  [19, 2, \"24-6=\"],\n  [19, 3, \"23+57=\"],\n  [19, 4, \"82-50=\"],\n];\n\nfor (const [row, col, newText] of updates) {\n  table.getCell(row, col).value = newText;\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rows = @(1, 1, 1, 1, 1, 2, 2, 2, 2, 2, 3, 3, 3, 3, 3, 4, 4, 4, 4, 4, 5, 5, 5, 5, 5, 6, 6, 6, 6, 6, 7, 7, 7, 7, 7, 8, 8, 8, 8, 8, 9, 9, 9, 9, 9, 10, 10, 10, 10, 10, 11, 11, 11, 11, 11, 12, 12, 12, 12, 12, 13, 13, 13, 13, 13, 14, 14, 14, 14, 14, 15, 15, 15, 15, 15, 16, 16, 16, 16, 16, 17, 17, 17, 17, 17, 18, 18, 18, 18, 18, 19, 19, 19, 19, 19, 20, 20, 20, 20, 20)\n$cols = @(1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5, 1, 2, 3, 4, 5)\n$texts = @(\"30+28=\", \"80-29=\", \"55-33=\", \"14+73=\", \"29-13=\", \"30+54=\", \"91-43=\", \"45-3=\", \"76-24=\", \"98-39=\", \"46+27=\", \"59+36=\", \"67-53=\", \"51-27=\", \"5+53=\", \"38+46=\", \"56-15=\", \"24+62=\", \"77-32=\", \"93-14=\", \"50-25=\", \"96-57=\", \"82-61=\", \"84+1=\", \"89-8=\", \"64-5=\", \"98-26=\", \"44+41=\", \"89-77=\", \"53+36=\", \"73-26=\", \"20+70=\", \"45+29=\", \"68-38=\", \"77-19=\", \"35+49=\", \"11+33=\", \"2+9=\", \"48-47=\", \"12+20=\", \"90-16=\", \"61-34=\", \"84+6=\", \"23-10=\", \"11+78=\", \"49-4=\", \"42+40=\", \"84-30=\", \"96-9=\", \"29+31=\", \"49-13=\", \"52+15=\", \"67+21=\", \"25+70=\", \"43-13=\", \"25+62=\", \"86-13=\", \"52+28=\", \"77-26=\", \"77-32=\", \"14+48=\", \"21-15=\", \"6+1=\", \"10+26=\", \"45-14=\", \"70-25=\", \"2+8=\", \"44-42=\", \"10+12=\", \"91-16=\", \"48-9=\", \"14+1=\", \"87-36=\", \"85-79=\", \"34-25=\", \"99-5=\", \"19+14=\", \"47-42=\", \"23+15=\", \"13+38=\", \"66-2=\", \"12-4=\", \"91-40=\", \"17+75=\", \"87-87=\", \"10+45=\", \"90-60=\", \"1+68=\", \"58-31=\", \"99-22=\", \"49-12=\", \"64-61=\", \"61-50=\", \"44-35=\", \"35-13=\", \"82+15=\", \"79-1=\", \"24-6=\", \"23+57=\", \"82-50=\")\n\nfor ($i = 0; $i -lt $rows.Length; $i++) {\n  $t.Cell($rows[$i], $cols[$i]).Range.Text = $texts[$i]\n}"}
